$d = $word.ActiveDocument

# NOTE: Paragraph/Range .Text in this host includes the trailing paragraph
# mark (Chr(13)), just like real Word - always trim it before comparing.
$cr = [char]13

# --- locate the paragraph that currently ends with the _GoBack bookmark ---
$src = $d.Paragraphs.Item(28)
if (($src.Range.Text.TrimEnd($cr)) -notlike "I did not make any diagrams*") {
    throw ("unexpected paragraph 28 text: " + $src.Range.Text)
}

# Remove the _GoBack bookmark from its current spot; we will recreate it
# later at its new home (right after the new "Problem:" paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- find "Socks in the dark" and the blank paragraph that follows it ---
$socks = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd($cr) -eq "Socks in the dark") {
        $socks = $i
        break
    }
}
if (-not $socks) { throw "could not find 'Socks in the dark' paragraph" }

$blankAfterSocks = $d.Paragraphs.Item($socks + 1)
if ($blankAfterSocks.Range.Text.TrimEnd($cr) -ne "") {
    throw "expected blank paragraph after 'Socks in the dark'"
}

# --- insert the new "Problem: ..." paragraph right after that blank line ---
$blankAfterSocks.Range.InsertParagraphAfter()
$problemPara = $d.Paragraphs.Item($socks + 2)
$problemPara.Range.Text = "Problem: You are in the dark, and you have a known number of socks with a known number of colors. How many socks would have to be picked to get a pair of any color, then how many would be needed to get a pair of all 3 colors."

# --- insert a fresh empty paragraph right after the "Problem:" paragraph; ---
# --- this paragraph will become the new home of the _GoBack bookmark     ---
$problemPara.Range.InsertParagraphAfter()
$bmPara = $d.Paragraphs.Item($socks + 3)

# Re-create the _GoBack bookmark collapsed at the start of that (still
# empty) paragraph. Bookmarks.Add on a truly empty range is unreliable in
# this host, so anchor it on a throwaway character and delete the
# character afterwards - the bookmark collapses to the right spot, same
# as Word does when you type + immediately delete at a bookmark.
$anchor = $bmPara.Range.Start
$tmp = $d.Range($anchor, $anchor)
$tmp.InsertAfter("x")
$tmpRange = $d.Range($anchor, $anchor + 1)
$d.Bookmarks.Add("_GoBack", $tmpRange)
$d.Range($anchor, $anchor + 1).Text = ""
